$wb = $excel.ActiveWorkbook
$wsHelp = $wb.Worksheets.Item("Help")
$wsHelp.Unprotect()
$wsHelp.Rows.Item(43).Delete()
$wsHelp.Protect($null, $true, $true, $true)
Write-Host "ProtectContents:" $wsHelp.ProtectContents
Write-Host "ProtectDrawingObjects:" $wsHelp.ProtectDrawingObjects
Write-Host "ProtectScenarios:" $wsHelp.ProtectScenarios
